$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 4622
$wsExpo.Range("F4").Value = 596
$wsExpo.Range("F6").Value = 1795
$wsExpo.Range("F8").Value = 720
$wsExpo.Range("F15").Value = 754
$wsExpo.Range("F16").Value = 544
$wsExpo.Range("F17").Value = 507
$wsExpo.Range("F19").Value = 155
$wsExpo.Range("F20").Value = 1549
$wsExpo.Range("F21").Value = 1184
$wsExpo.Range("F22").Value = 385
$wsExpo.Range("F23").Value = 2496
$wsExpo.Range("F25").Value = 1530
$wsExpo.Range("F29").Value = 4209

# --- Sheet "演出" (Performance) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F17").Value = 279

# --- Sheet "本地生活" (Local life) ---
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F7").Value = 226

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 226
$wsAll.Range("F8").Value = 4622
$wsAll.Range("F9").Value = 596
$wsAll.Range("F11").Value = 1795
$wsAll.Range("F12").Value = 720

# Row 19 in "全部类型" is replaced wholesale: the old
# "热血之巅·突破次元壁" entry is swapped for the "天空之城" one.
$wsAll.Range("C19").Value = "上海·《天空之城》宫崎骏&久石让经典作品动漫视听音乐会"
$wsAll.Range("E19").Value = "2024.03.31 10:30-03.31 12:00"
$wsAll.Range("F19").Value = 25
$wsAll.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=81660"
$wsAll.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202402/QKmfdsEM1706853934802.jpeg"

$wsAll.Range("F22").Value = 755
$wsAll.Range("F23").Value = 544
$wsAll.Range("F24").Value = 507
$wsAll.Range("F26").Value = 155
$wsAll.Range("F28").Value = 279
$wsAll.Range("F32").Value = 1549
$wsAll.Range("F33").Value = 1184
$wsAll.Range("F34").Value = 385
$wsAll.Range("F37").Value = 2496
$wsAll.Range("F43").Value = 1530
$wsAll.Range("F48").Value = 4209
